$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 31; $r++) {
    $cell = $ws.Cells.Item($r, 58)
    $cell.Value = "'2014-04-25"
}
